$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Force D:E data range to Text format before writing values so that
# numeric-looking strings (e.g. "530.71") are preserved as text, matching
# the original inlineStr cell type instead of being auto-converted to numbers.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = '58.611.06'
$ws.Range("E2").Value = '  +1.77%  '
$ws.Range("D3").Value = '3.149.50'
$ws.Range("E3").Value = '  +1.58%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '530.71'
$ws.Range("E5").Value = '  +0.50%  '
$ws.Range("D6").Value = '140.02'
$ws.Range("E6").Value = '  +1.93%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("D8").Value = '0.518'
$ws.Range("E8").Value = '  +10.91%  '
$ws.Range("E9").Value = '  +1.03%  '
$ws.Range("E10").Value = '  +3.01%  '
$ws.Range("D11").Value = '0.422'
$ws.Range("E11").Value = '  +3.68%  '
$ws.Range("E12").Value = '  +2.24%  '
$ws.Range("D13").Value = '3.692.71'
$ws.Range("E13").Value = '  +1.52%  '
$ws.Range("D14").Value = '25.77'
$ws.Range("E14").Value = '  +1.49%  '
$ws.Range("E15").Value = '  +5.46%  '
$ws.Range("D16").Value = '58.645.59'
$ws.Range("E16").Value = '  +1.73%  '
$ws.Range("D17").Value = '3.152.43'
$ws.Range("E17").Value = '  +2.13%  '
$ws.Range("E18").Value = '  +4.60%  '
$ws.Range("E19").Value = '  +3.44%  '
$ws.Range("E20").Value = '  +3.00%  '
$ws.Range("D21").Value = '371.77'
$ws.Range("E21").Value = '  +6.16%  '
$ws.Range("D22").Value = '5.80'
$ws.Range("E22").Value = '  +1.89%  '
$ws.Range("D23").Value = '1.00'
$ws.Range("D24").Value = '69.61'
$ws.Range("E24").Value = '  +1.04%  '
$ws.Range("E25").Value = '  +2.08%  '
$ws.Range("E26").Value = '  +0.27%  '
$ws.Range("E27").Value = '  +0.20%  '
$ws.Range("D28").Value = '7.98'
$ws.Range("E28").Value = '  +12.00%  '
$ws.Range("E29").Value = '  +0.50%  '
$ws.Range("D30").Value = '6.13'
$ws.Range("E30").Value = '  +2.56%  '
$ws.Range("E31").Value = '  +0.97%  '
$ws.Range("D32").Value = '21.89'
$ws.Range("E32").Value = '  +3.70%  '
$ws.Range("E33").Value = '  +6.19%  '
$ws.Range("E34").Value = '  +2.48%  '
$ws.Range("D35").Value = '159.85'
$ws.Range("E35").Value = '  +0.59%  '
$ws.Range("D36").Value = '6.26'
$ws.Range("E36").Value = '  +3.97%  '
$ws.Range("D37").Value = '1.36'
$ws.Range("E37").Value = '  +9.29%  '
$ws.Range("D38").Value = '25.19'
$ws.Range("E38").Value = '  -1.96%  '
$ws.Range("D39").Value = '2.655.77'
$ws.Range("E39").Value = '  +10.75%  '
$ws.Range("D40").Value = '1.67'
$ws.Range("E40").Value = '  +3.16%  '
$ws.Range("D41").Value = '0.0686'
$ws.Range("E41").Value = '  +3.16%  '
$ws.Range("E42").Value = '  +3.97%  '
$ws.Range("D43").Value = '0.709'
$ws.Range("E43").Value = '  +2.12%  '
$ws.Range("E44").Value = '  +9.00%  '
$ws.Range("D45").Value = '38.41'
$ws.Range("E45").Value = '  +4.01%  '
$ws.Range("E46").Value = '  -0.10%  '
$ws.Range("D47").Value = '3.191.20'
$ws.Range("E47").Value = '  +1.48%  '
$ws.Range("D48").Value = '0.104'
$ws.Range("E48").Value = '  +13.65%  '
$ws.Range("D49").Value = '0.979'
$ws.Range("E49").Value = '  +2.82%  '
$ws.Range("E50").Value = '  +2.57%  '
$ws.Range("D51").Value = '20.19'
$ws.Range("E51").Value = '  +3.34%  '

# Restore the default cell style so the written cells match the original
# workbook formatting (no explicit style index / number format).
$dataRange.Style = "Normal"
